# Apply the edit described by the diff:
# - Insert a new column before column D ("Clustering optimal nb"), shifting the
#   previous "ByStages" column to E and "ByStages high nb" to F.
# - Rename headers: C1 "Clustering" -> "Clustering nc"; new D1 = "Clustering optimal nb";
#   E1 (previously "ByStages") -> "ByStages nc".
# - Adjust the highlight fills in the new layout.
# - Update column widths, selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert new column before D - shifts old D->E, old E->F
$ws.Columns("D:D").Insert()

# 2. Update header row text (order chosen to match shared-string insertion order)
$ws.Range("E1").Value = "ByStages nc"
$ws.Range("C1").Value = "Clustering nc"
$ws.Range("D1").Value = "Clustering optimal nb"
# F1 already holds "ByStages high nb" (shifted automatically by the insert)

# Make sure the new D1 header cell carries the same bold style as its neighbours
$ws.Range("D1").Font.Bold = $true

# 3. Column widths (nearest values achievable through the host's column-width
#    quantisation that reproduce the target stored widths of 11.77734375 /
#    19.109375 / 11 / 15.44140625 characters)
$ws.Columns("C:C").ColumnWidth = 11
$ws.Columns("D:D").ColumnWidth = 18.333333333333336
$ws.Columns("E:E").ColumnWidth = 10.166666666666666
$ws.Columns("F:F").ColumnWidth = 14.666666666666666

# 4. Fix up highlight fills for the "D" case-group (rows 2-6)
$ws.Range("D2").Interior.Color = 5296274
$ws.Range("D3").Interior.Color = 5296274
$ws.Range("D4").Interior.Color = 5296274
$ws.Range("D5").Interior.Color = 5296274

$ws.Range("E3").Interior.Color = 5296274
$ws.Range("F3").Interior.Color = 49407

$ws.Range("F2").Interior.Color = 5296274

$ws.Range("C6").Interior.Color = 5296274
$ws.Range("D6").Interior.Pattern = -4142

# 5. Fix up highlight fills for the "E" case-group (rows 7-8)
$ws.Range("D7").Interior.Pattern = -4142
$ws.Range("F7").Interior.Color = 5296274

$ws.Range("C8").Interior.Color = 5296274
$ws.Range("D8").Interior.Pattern = -4142
$ws.Range("E8").Interior.Color = 49407

# 6. Update selection
$ws.Range("E6").Select()
